{"js": "// Portuguese (pt) translation pass for the \"documents verified\" partner\n// email template. Each English source string is located with\n// Body.search (exact, case-sensitive, literal text) and its single\n// result range has its text swapped for the Portuguese copy via\n// insertText(..., Word.InsertLocation.replace), which rewrites the text\n// of the existing run(s) in place without touching surrounding runs,\n// formatting or other document structure.\n\nasync function replaceOnce(body, needle, replacement, occurrenceIndex) {\n  const idx = occurrenceIndex || 0;\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length <= idx) {\n    throw new Error(\n      \"Expected occurrence \" + idx + \" of \" + JSON.stringify(needle) +\n      \" but only found \" + results.items.length + \" match(es).\"\n    );\n  }\n  results.items[idx].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) Standalone \"English\" heading above the table (2nd \"English\" in the\n//    doc; the 1st is inside the language-switcher hyperlink and must stay).\nawait replaceOnce(body, \"English\", \"Ingl\u00eas\", 1);\n\n// 2) Big centered headline.\nawait replaceOnce(\n  body,\n  \"Your documents have been verified!\",\n  \"Os seus documentos foram verificados!\"\n);\n\n// 3) Greeting line.\nawait replaceOnce(body, \"Hi \", \"Ol\u00e1 \");\nawait replaceOnce(body, \"[PARTNER NAME]\", \"[NOME DO PARCEIRO]\");\n\n// 4) \"We've reviewed...\" paragraph (split across three runs around the\n//    highlighted [EVENT NAME] placeholder).\nawait replaceOnce(\n  body,\n  \"We\u2019ve reviewed the documents you\u2019ve sent us for the \",\n  \"Analis\u00e1mos os documentos que nos enviou para a \"\n);\n// 2nd \"[EVENT NAME]\" occurrence (1st is in the Subject line paragraph).\nawait replaceOnce(body, \"[EVENT NAME]\", \"[NOME DO EVENTO]\", 1);\nawait replaceOnce(\n  body,\n  \" and all of them have been verified! \",\n  \" e foram todos verificados! \"\n);\n\n// 5) \"We'll be sending...\" paragraph.\nawait replaceOnce(\n  body,\n  \"We\u2019ll be sending out more details about the event soon, including the agenda and travel arrangements, so make sure to check your inbox regularly.\",\n  \"Em breve, enviaremos mais pormenores sobre o evento, incluindo a agenda e os preparativos para a viagem, por isso n\u00e3o se esque\u00e7a de consultar regularmente o seu e-mail.\"\n);\n\n// 6) \"If you have any questions, please contact us via live chat or\n//    WhatsApp.\" paragraph.\nawait replaceOnce(\n  body,\n  \"If you have any questions, please contact us via \",\n  \"Para mais informa\u00e7\u00f5es, contacte-nos atrav\u00e9s de \"\n);\n// 1st \" or \" occurrence (between the live-chat and WhatsApp hyperlinks).\nawait replaceOnce(body, \" or \", \" ou \", 0);\n\n// 7) \"If you have any questions, please contact your country manager...\"\n//    paragraph.\nawait replaceOnce(\n  body,\n  \"If you have any questions, please contact your country manager, \",\n  \"Para mais quest\u00f5es, pode tamb\u00e9m contactar o seus gestor de parcerias \"\n);\nawait replaceOnce(body, \", at \", \", em \");\n// After the two replacements above, the only remaining \" or \" is the one\n// between [EMAIL ADDRESS] and [WHATSAPP NO].\nawait replaceOnce(body, \" or \", \" ou \", 0);\n", "ps1": "# Portuguese (pt) translation pass for the \"documents verified\" partner\n# email template.\n#\n# For each English source string we walk $d.Content with Find.Execute\n# (literal text, not wildcards) to the requested occurrence (1-based,\n# matching the order the strings appear in the document), append the\n# Portuguese replacement right after the found range with InsertAfter,\n# then delete the original (now-stale) range.\n#\n# We deliberately avoid just setting `$range.Text = $ReplaceText`: when\n# the found range starts exactly at a comment anchor (w:commentRangeStart)\n# or right after a differently-formatted run (e.g. a hyperlink), a plain\n# Range.Text assignment here re-anchors that marker/formatting to the\n# wrong side of the new text. Insert-after-then-delete-the-old-range\n# leaves every other marker/run exactly where it was.\n$d = $word.ActiveDocument\n\nfunction Replace-NthOccurrence {\n    param(\n        $Doc,\n        [string]$FindText,\n        [string]$ReplaceText,\n        [int]$Occurrence = 1\n    )\n\n    $range = $Doc.Content\n    for ($i = 0; $i -lt $Occurrence; $i++) {\n        $found = $range.Find.Execute($FindText, $true)\n        if (-not $found) {\n            throw \"Replace-NthOccurrence: could not find occurrence $($i+1) of '$FindText'\"\n        }\n    }\n\n    $origStart = $range.Start\n    $origEnd = $range.End\n    $range.InsertAfter($ReplaceText)\n    $old = $Doc.Range($origStart, $origEnd)\n    $old.Delete()\n}\n\n# 1) Standalone \"English\" heading above the table (2nd \"English\" in the\n#    doc; the 1st is inside the language-switcher hyperlink and must stay).\nReplace-NthOccurrence $d \"English\" \"Ingl\u00eas\" 2\n\n# 2) Big centered headline.\nReplace-NthOccurrence $d \"Your documents have been verified!\" \"Os seus documentos foram verificados!\" 1\n\n# 3) Greeting line.\nReplace-NthOccurrence $d \"Hi \" \"Ol\u00e1 \" 1\nReplace-NthOccurrence $d \"[PARTNER NAME]\" \"[NOME DO PARCEIRO]\" 1\n\n# 4) \"We've reviewed...\" paragraph (split across three runs around the\n#    highlighted [EVENT NAME] placeholder).\nReplace-NthOccurrence $d \"We\u2019ve reviewed the documents you\u2019ve sent us for the \" \"Analis\u00e1mos os documentos que nos enviou para a \" 1\n# 2nd \"[EVENT NAME]\" occurrence (1st is in the Subject line paragraph).\nReplace-NthOccurrence $d \"[EVENT NAME]\" \"[NOME DO EVENTO]\" 2\nReplace-NthOccurrence $d \" and all of them have been verified! \" \" e foram todos verificados! \" 1\n\n# 5) \"We'll be sending...\" paragraph.\nReplace-NthOccurrence $d \"We\u2019ll be sending out more details about the event soon, including the agenda and travel arrangements, so make sure to check your inbox regularly.\" \"Em breve, enviaremos mais pormenores sobre o evento, incluindo a agenda e os preparativos para a viagem, por isso n\u00e3o se esque\u00e7a de consultar regularmente o seu e-mail.\" 1\n\n# 6) \"If you have any questions, please contact us via live chat or\n#    WhatsApp.\" paragraph.\nReplace-NthOccurrence $d \"If you have any questions, please contact us via \" \"Para mais informa\u00e7\u00f5es, contacte-nos atrav\u00e9s de \" 1\n# 1st \" or \" occurrence (between the live-chat and WhatsApp hyperlinks).\nReplace-NthOccurrence $d \" or \" \" ou \" 1\n\n# 7) \"If you have any questions, please contact your country manager...\"\n#    paragraph.\nReplace-NthOccurrence $d \"If you have any questions, please contact your country manager, \" \"Para mais quest\u00f5es, pode tamb\u00e9m contactar o seus gestor de parcerias \" 1\nReplace-NthOccurrence $d \", at \" \", em \" 1\n# After the two replacements above, the only remaining \" or \" is the one\n# between [EMAIL ADDRESS] and [WHATSAPP NO].\nReplace-NthOccurrence $d \" or \" \" ou \" 1\n"}
